# Updates the crypto price/volume table to the latest scrape snapshot.
# For each changed cell we force text formatting (NumberFormat "@") before
# assigning the value so numeric-looking strings (e.g. "409.88") are not
# auto-converted to Excel numbers, matching the source data's inline-string
# cell type. The style is reset to "Normal" afterwards so the NumberFormat
# tweak does not leave a stray cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '62.057.91'
Set-TextValue 'E2' '  +1.01%  '

# Row 3
Set-TextValue 'D3' '3.421.33'
Set-TextValue 'E3' '  +1.21%  '

# Row 4
Set-TextValue 'E4' '  +0.01%  '

# Row 5
Set-TextValue 'D5' '409.88'
Set-TextValue 'E5' '  +1.05%  '

# Row 6
Set-TextValue 'D6' '128.19'
Set-TextValue 'E6' '  -4.49%  '

# Row 7
Set-TextValue 'D7' '0.626'
Set-TextValue 'E7' '  +5.79%  '

# Row 8
Set-TextValue 'E8' '  -0.08%  '

# Row 9
Set-TextValue 'D9' '0.747'
Set-TextValue 'E9' '  +11.21%  '

# Row 10
Set-TextValue 'E10' '  +18.40%  '

# Row 11
Set-TextValue 'D11' '43.02'
Set-TextValue 'E11' '  +0.56%  '

# Row 12
Set-TextValue 'E12' '  -0.46%  '

# Row 13
Set-TextValue 'D13' '3.952.68'
Set-TextValue 'E13' '  +1.44%  '

# Row 14
Set-TextValue 'D14' '21.13'
Set-TextValue 'E14' '  +6.99%  '

# Row 15
Set-TextValue 'D15' '8.90'
Set-TextValue 'E15' '  +5.75%  '

# Row 16
Set-TextValue 'D16' '0.0000207'
Set-TextValue 'E16' '  +62.87%  '

# Row 17
Set-TextValue 'D17' '3.410.28'
Set-TextValue 'E17' '  +0.87%  '

# Row 18
Set-TextValue 'D18' '12.60'
Set-TextValue 'E18' '  +14.70%  '

# Row 19
Set-TextValue 'D19' '1.08'
Set-TextValue 'E19' '  +5.11%  '

# Row 20
Set-TextValue 'D20' '62.099.44'
Set-TextValue 'E20' '  +1.21%  '

# Row 21
Set-TextValue 'D21' '405.17'
Set-TextValue 'E21' '  +28.58%  '

# Row 22
Set-TextValue 'D22' '90.54'
Set-TextValue 'E22' '  +6.52%  '

# Row 23
Set-TextValue 'E23' '  -0.49%  '

# Row 24
Set-TextValue 'D24' '13.56'
Set-TextValue 'E24' '  +5.44%  '

# Row 25
Set-TextValue 'E25' '  +2.76%  '

# Row 26
Set-TextValue 'D26' '32.98'
Set-TextValue 'E26' '  +11.70%  '

# Row 27
Set-TextValue 'B27' 'Filecoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D27' '8.59'
Set-TextValue 'E27' '  +3.19%  '

# Row 28
Set-TextValue 'B28' 'LEO'
Set-TextValue 'C28' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D28' '4.80'
Set-TextValue 'E28' '  +0.52%  '

# Row 29
Set-TextValue 'D29' '7.64'
Set-TextValue 'E29' '  +0.69%  '

# Row 30
Set-TextValue 'E30' '  +1.89%  '

# Row 31
Set-TextValue 'E31' '  +1.27%  '

# Row 32
Set-TextValue 'B32' 'Kaspa'
Set-TextValue 'C32' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D32' '0.172'
Set-TextValue 'E32' '  +0.65%  '

# Row 33
Set-TextValue 'B33' 'InjectiveProtocol'
Set-TextValue 'C33' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D33' '43.89'
Set-TextValue 'E33' '  +6.57%  '

# Row 34
Set-TextValue 'D34' '11.88'
Set-TextValue 'E34' '  +4.49%  '

# Row 35
Set-TextValue 'E35' '  -0.03%  '

# Row 36
Set-TextValue 'D36' '0.0499'
Set-TextValue 'E36' '  +3.72%  '

# Row 37
Set-TextValue 'D37' '52.98'
Set-TextValue 'E37' '  +2.51%  '

# Row 38
Set-TextValue 'E38' '  +0.14%  '

# Row 39
Set-TextValue 'E39' '  -0.92%  '

# Row 40
Set-TextValue 'E40' '  -0.89%  '

# Row 41
Set-TextValue 'D41' '0.131'
Set-TextValue 'E41' '  +5.91%  '

# Row 42
Set-TextValue 'D42' '0.314'
Set-TextValue 'E42' '  +5.93%  '

# Row 43
Set-TextValue 'D43' '141.19'
Set-TextValue 'E43' '  +0.83%  '

# Row 44
Set-TextValue 'D44' '1.98'
Set-TextValue 'E44' '  -0.07%  '

# Row 45
Set-TextValue 'D45' '4.02'
Set-TextValue 'E45' '  -0.09%  '

# Row 46
Set-TextValue 'D46' '2.38'
Set-TextValue 'E46' '  +6.52%  '

# Row 47
Set-TextValue 'D47' '16.77'
Set-TextValue 'E47' '  +0.73%  '

# Row 48
Set-TextValue 'D48' '21.96'
Set-TextValue 'E48' '  +2.99%  '

# Row 49
Set-TextValue 'D49' '2.113.39'
Set-TextValue 'E49' '  -0.13%  '

# Row 50
Set-TextValue 'B50' 'ThetaToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D50' '1.94'
Set-TextValue 'E50' '  +1.23%  '

# Row 51
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.127'
Set-TextValue 'E51' '  +13.96%  '
